# Updated "Tasks 01-28 to 02-04" task list sheet:
#  - logged Time Spent for a couple of in-progress tasks and marked them Done
#  - logged Time Spent for two more in-progress tasks (left as TODO)
#  - added a new task row for "Convert unequippedItems to a List"
#  - refreshed the sheet's active cell selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 ("Filter InventoryItems based on itemType for Inventory Scene") ---
#     Time Spent = 2, Status moves from TODO -> Done
$ws.Range("C6").Value = 2
$ws.Range("F2").Copy()
$ws.Range("F6").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F6").Value = "Done"

# --- Row 7 ("Select inventoryItem with touch in Inventory Scene") ---
#     Time Spent = 1, Status moves from TODO -> Done
$ws.Range("C7").Value = 1
$ws.Range("F2").Copy()
$ws.Range("F7").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F7").Value = "Done"

# --- Row 8 ("Create ComparedItem Game Object") : Time Spent = 2, status unchanged ---
$ws.Range("C8").Value = 2

# --- Row 9 ("Create SelectedItem Game Object") : Time Spent = 1, status unchanged ---
$ws.Range("C9").Value = 1

# --- New row 13: "Convert unequippedItems to a List" ---
$ws.Range("A13").Value = "Convert unequippedItems to a List"
$ws.Range("B13").Value = 1
$ws.Range("E13").Value = "James"

$ws.Range("F8").Copy()
$ws.Range("F13").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F13").Value = "TODO"

# --- Refresh selection shown when the sheet is opened ---
$ws.Range("C8").Select()

$wb.Save()
